$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The new experiment (run 61, "Overtain test on small complex capsunet")
# is logged as a 4-sub-row block (0_voice / 0_background / 1_voice /
# 1_background), the same shape as the existing run-52 block in rows
# 15-18. Seed the formatting for rows 21-24 by copying each source row's
# formats+values cell-by-cell (keeps the style table from growing), then
# overwrite with the real values for the new run.
$ws.Range("A15").Copy($ws.Range("A21"))
$ws.Range("B15").Copy($ws.Range("B21"))
$ws.Range("C15:H15").Copy($ws.Range("C21:H21"))
$ws.Range("I15").Copy($ws.Range("I21"))

$ws.Range("A16").Copy($ws.Range("A22"))
$ws.Range("B16").Copy($ws.Range("B22"))
$ws.Range("C16:H16").Copy($ws.Range("C22:H22"))
$ws.Range("I16").Copy($ws.Range("I22"))

$ws.Range("A17").Copy($ws.Range("A23"))
$ws.Range("B17").Copy($ws.Range("B23"))
$ws.Range("C17:H17").Copy($ws.Range("C23:H23"))
$ws.Range("I17").Copy($ws.Range("I23"))

$ws.Range("A18").Copy($ws.Range("A24"))
$ws.Range("B18").Copy($ws.Range("B24"))
$ws.Range("C18:H18").Copy($ws.Range("C24:H24"))
$ws.Range("I18").Copy($ws.Range("I24"))

# Row 21 - "0_voice" sub-row (header of the new experiment block)
$ws.Range("A21").Value = 43373
$ws.Range("B21").Value = 61
$ws.Range("C21").Value = "0_voice"
$ws.Range("D21").Value = 0.147776272892951
$ws.Range("E21").Value = 5.62416140137290022
$ws.Range("F21").Value = 5.67600894251940957
$ws.Range("G21").Value = 26.41188428278550049
$ws.Range("H21").Value = 0.00859465250988306
$ws.Range("I21").Value = "Overtain test on small complex capsunet"

# Row 22 - "0_background" sub-row
$ws.Range("C22").Value = "0_background"
$ws.Range("D22").Value = 0.147776272892951
$ws.Range("E22").Value = 14.03600391987320073
$ws.Range("F22").Value = 14.31212843302770032
$ws.Range("G22").Value = 26.44145597086199828
$ws.Range("H22").Value = 0.0368738486030178

# Row 23 - "1_voice" sub-row
$ws.Range("C23").Value = "1_voice"
$ws.Range("D23").Value = 0.12048309892416
$ws.Range("E23").Value = 6.4002952428711497
$ws.Range("F23").Value = 6.96789942232423964
$ws.Range("G23").Value = 17.12628030908319943
$ws.Range("H23").Value = 0.78472849400813605

# Row 24 - "1_background" sub-row
$ws.Range("C24").Value = "1_background"
$ws.Range("D24").Value = 0.12048309892416
$ws.Range("E24").Value = 14.01941207764300046
$ws.Range("F24").Value = 14.66750488069880021
$ws.Range("G24").Value = 22.94366919112480119
$ws.Range("H24").Value = 0.0202820063727339

# Re-create the merges for the new block (A/B/I columns span the 4 sub-rows).
$ws.Range("A21:A24").Merge()
$ws.Range("B21:B24").Merge()
$ws.Range("I21:I24").Merge()

# Move the active selection to mimic where the author's cursor ended up.
$ws.Range("I31").Select()
